$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13: F13 formula changes from =12*2 to =12, and a new constant 12 is placed in L13
$ws.Range("F13").Formula = "=12"
$ws.Range("L13").Value = 12

# Row 29: the value 20 moves from F29 to L29
$ws.Range("F29").ClearContents()
$ws.Range("L29").Value = 20

# Update selection to L32 to match the saved cursor position
$ws.Range("L32").Select()
